$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header in J1: CO2_constraint -> Global_CO2_constraint
$ws.Range("J1").Value = "Global_CO2_constraint"

# Add new column header N1
$ws.Range("N1").Value = "2050CO2"

# Fill N2:N7 with 0 (existing rows get the new column populated with 0)
$ws.Range("N2:N7").Value = 0

# Add new row 7 data
$ws.Range("B7").Value = "NLP"
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = "Ipopt"
$ws.Range("F7").Value = 20451
$ws.Range("G7").Value = 1318
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.01
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0

# Add new row 8 data (partial row)
$ws.Range("B8").Value = "NLP"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = "Ipopt"
$ws.Range("N8").Value = 1

# Add new rows 9-11 (only E filled)
$ws.Range("E9").Value = "Ipopt"
$ws.Range("E10").Value = "Ipopt"
$ws.Range("E11").Value = "Ipopt"

# Recompute best-fit width for column J (header text got longer, bestFit widens it)
$ws.Columns.Item(10).ColumnWidth = 19.3307291667

# Update selection to F8
$ws.Range("F8").Select()
